$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# K114 = 0 (new cell in existing row)
$ws.Range("K114").Value = 0

# Row 115: DQ Tau
$ws.Range("A115").Value = "DQ Tau"
$ws.Range("A115").Font.Color = 0x010101
$ws.Range("B115").Value = "J04465305+1700001"
$ws.Range("B115").Font.Color = 0x010101
$ws.Range("C115").Value = "HBC 72"
$ws.Range("D115").Value = "04h46m53.06s"
$ws.Range("D115").Font.Color = 0x010101
$ws.Range("E115").Value = "'+17d00m00.14s"
$ws.Range("E115").Font.Color = 0x010101
$ws.Range("F115").Value = "Taurus"
$ws.Range("F115").Font.Color = 0x010101
$ws.Range("G115").Value = 196
$ws.Range("H115").Value = "M0-1"
$ws.Range("I115").Value = 1.21
$ws.Range("K115").Value = 0.5
$ws.Range("N115").Value = 2
$ws.Range("O115").Value = 0
$ws.Range("P115").Value = 0
$ws.Range("Q115").Value = 0
$ws.Range("R115").Value = 0
$ws.Range("S115").Value = 0
$ws.Range("T115").Value = 0
$ws.Range("U115").Value = 0
$ws.Range("V115").Value = 0
$ws.Range("W115").Value = 12161
$ws.Range("X115").Value = "COS/G160M-COS/G230L"
$ws.Range("Y115").Value = 6
$ws.Range("Z115").Value = "czekala+2016"

# Row 116: KK Oph
$ws.Range("A116").Value = "KK Oph"
$ws.Range("A116").Font.Color = 0x010101
$ws.Range("B116").Value = "J17100811-2715190"
$ws.Range("B116").Font.Color = 0x010101
$ws.Range("C116").Value = "HBC 273"
$ws.Range("D116").Value = "17h10m08.12s"
$ws.Range("D116").Font.Color = 0x010101
$ws.Range("E116").Value = "'-27d15m18.80s"
$ws.Range("E116").Font.Color = 0x010101
$ws.Range("F116").Value = "other"
$ws.Range("F116").Font.Color = 0x010101
$ws.Range("G116").Value = 160
$ws.Range("H116").Value = "A6V + G6V"
$ws.Range("I116").Value = 1.6
$ws.Range("K116").Value = 1.6
$ws.Range("L116").Value = 12.86
$ws.Range("M116").Value = 12.22
$ws.Range("N116").Value = 2
$ws.Range("O116").Value = 0
$ws.Range("P116").Value = 0
$ws.Range("Q116").Value = 0
$ws.Range("R116").Value = 0
$ws.Range("S116").Value = 0
$ws.Range("T116").Value = 0
$ws.Range("U116").Value = 0
$ws.Range("V116").Value = 0
$ws.Range("W116").Value = 12996
$ws.Range("X116").Value = "COS/G130M-COS/G160M-COS/G230L"
$ws.Range("Y116").Value = 6
$ws.Range("Z116").Value = "carmona2007"
